$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Remove the old placeholder hyperlinks (PruebaN@pruebmail.com)
#    so we can rebuild them pointing at the real e-mail addresses.
# ------------------------------------------------------------------
$ws.Hyperlinks.Delete()

# ------------------------------------------------------------------
# 2. Update the e-mail text shown in column L for the existing rows
#    (2-6), replacing the test placeholders with the real addresses.
# ------------------------------------------------------------------
$ws.Range("L2").Value = "quintana.nicolas@javeriana.edi.co"
$ws.Range("L3").Value = "jecheverry@javeriana.edu.co"
$ws.Range("L4").Value = "Prueba3@pruebmail.com"
$ws.Range("L5").Value = "valentina.lopez@javeriana.edu.co"
$ws.Range("L6").Value = "mateo.fernandez@.hotmail.com"

# ------------------------------------------------------------------
# 3. Add row 7 - a duplicate of row 6's survey answers for a new
#    tutor/tutee pair, with its own e-mail address in column L.
# ------------------------------------------------------------------
$ws.Range("A7").Value = "2024-01-07 07:52:23"
$ws.Range("B7").Value = "2024-01-17 04:44:29"
$ws.Range("C7").Value = "Completa"
$ws.Range("D7").Value = "192.168.1.3"
$ws.Range("E7").Value = 95
$ws.Range("F7").Value = 517
$ws.Range("G7").Value = "Sí"
$ws.Range("H7").Value = "2025-02-12 18:19:00"
$ws.Range("I7").Value = "RESP-4175"
$ws.Range("J7").Value = "Apellido19"
$ws.Range("K7").Value = "Nombre30"
$ws.Range("L7").Value = "camila.rodriguez@javeriana.edu.co"
$ws.Range("W7").Value = 5555
$ws.Range("AC7").Value = "Muy motivado"
$ws.Range("AD7").Value = "Malo"
$ws.Range("AM7").Value = "Insatisfecho"
$ws.Range("AN7").Value = "Problemas técnicos"
$ws.Range("AP7").Value = 5
$ws.Range("AQ7").Value = "Gracias"

# ------------------------------------------------------------------
# 4. Re-create the hyperlinks on column L (rows 2-7), one per e-mail.
# ------------------------------------------------------------------
$ws.Hyperlinks.Add($ws.Range("L2"), "mailto:quintana.nicolas@javeriana.edi.co")
$ws.Hyperlinks.Add($ws.Range("L3"), "mailto:jecheverry@javeriana.edu.co")
$ws.Hyperlinks.Add($ws.Range("L4"), "mailto:Prueba3@pruebmail.com")
$ws.Hyperlinks.Add($ws.Range("L5"), "mailto:valentina.lopez@javeriana.edu.co")
$ws.Hyperlinks.Add($ws.Range("L6"), "mailto:mateo.fernandez@.hotmail.com")
$ws.Hyperlinks.Add($ws.Range("L7"), "mailto:camila.rodriguez@javeriana.edu.co")

# Restore the original "Hipervínculo" cell style on the hyperlink
# cells (Hyperlinks.Add otherwise mints a brand-new duplicate style).
$ws.Range("L2:L7").Style = "Hipervínculo"

# ------------------------------------------------------------------
# 5. Update the view: drop the old top-left-cell pin and move the
#    selection to F13, matching the saved workbook state.
# ------------------------------------------------------------------
$ws.Range("F13").Select()
